$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1")
try {
  $fcs = $rng.FormatConditions
  Write-Host "FormatConditions count:" $fcs.Count
} catch { Write-Host "err" $_.Exception.Message }
